$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 20: 2018-08-17, python / Hankerank tasks on python (repeat of row 19)
$ws.Range("A20").Value = 43329
$ws.Range("A20").NumberFormat = "m/d/yy"
$ws.Range("B20").Value = "python"
$ws.Range("C20").Value = "Hankerank tasks on python"

# Row 21: 2018-08-18, python / word occurrence counting tasks
$ws.Range("A21").Value = 43330
$ws.Range("A21").NumberFormat = "m/d/yy"
$ws.Range("B21").Value = "python"
$ws.Range("C21").Value = "tasks:count word occurrences(case sensitive) and count word occurrences(case insensitive"

# Row 22: 2018-08-19, Python / extract links from webpage task
$ws.Range("A22").Value = 43331
$ws.Range("A22").NumberFormat = "m/d/yy"
$ws.Range("B22").Value = "Python"
$ws.Range("C22").Value = "tasks: Extract links from a webpage"

$ws.Range("C22").Select()
